$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.841.97"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "2.366.09"
$ws.Range("E3").Value = "  +2.21%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'301.28"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").Value = "'95.52"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("E7").Value = "  -0.67%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "'0.485"
$ws.Range("E9").Value = "  -1.68%  "
$ws.Range("D10").Value = "'33.93"
$ws.Range("E10").Value = "  -1.38%  "
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("E12").Value = "  +2.90%  "
$ws.Range("D13").Value = "'18.31"
$ws.Range("E13").Value = "  -3.56%  "
$ws.Range("D14").Value = "'6.73"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "2.734.98"
$ws.Range("E15").Value = "  +2.10%  "
$ws.Range("D16").Value = "2.368.76"
$ws.Range("E16").Value = "  +2.39%  "
$ws.Range("E17").Value = "  +1.23%  "
$ws.Range("D18").Value = "42.807.79"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").Value = "'12.07"
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("D20").Value = "'6.27"
$ws.Range("E20").Value = "  +1.89%  "
$ws.Range("E21").Value = "  -0.91%  "
$ws.Range("D22").Value = "'67.89"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "'234.83"
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("D24").Value = "'2.21"
$ws.Range("E24").Value = "  -1.96%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("D27").Value = "'24.71"
$ws.Range("E27").Value = "  +1.39%  "
$ws.Range("E28").Value = "  +0.38%  "
$ws.Range("D29").Value = "'9.21"
$ws.Range("E29").Value = "  +0.95%  "
$ws.Range("D30").Value = "'31.42"
$ws.Range("E30").Value = "  -2.53%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  +0.48%  "
$ws.Range("D33").Value = "'0.0731"
$ws.Range("E33").Value = "  +4.79%  "
$ws.Range("D34").Value = "'17.35"
$ws.Range("E34").Value = "  -3.25%  "
$ws.Range("E35").Value = "  +4.63%  "
$ws.Range("E36").Value = "  +4.39%  "
$ws.Range("D37").Value = "'4.34"
$ws.Range("E37").Value = "  -2.43%  "
$ws.Range("E38").Value = "  -1.57%  "
$ws.Range("E39").Value = "  +1.75%  "
$ws.Range("D40").Value = "'22.07"
$ws.Range("E40").Value = "  +5.82%  "
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("D42").Value = "'118.03"
$ws.Range("E42").Value = "  -28.95%  "
$ws.Range("D43").Value = "1.932.09"
$ws.Range("E43").Value = "  +0.26%  "
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("E45").Value = "  +1.81%  "
$ws.Range("D46").Value = "'2.72"
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("D47").Value = "'9.16"
$ws.Range("E47").Value = "  -9.58%  "
$ws.Range("D48").Value = "2.599.73"
$ws.Range("E48").Value = "  +2.14%  "
$ws.Range("E49").Value = "  +2.17%  "
$ws.Range("D50").Value = "'72.03"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").Value = "'51.77"
$ws.Range("E51").Value = "  -3.01%  "
